# Adding functional test case for registration email text field
$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (TC005) so it lands in
# the same tab position the author's "TC007" sheet occupies.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TC007"

# Header
$ws.Range("A1").Value = "Email"

# Valid email -> Excel auto-recognizes it and we wire up the mailto: link
$ws.Range("A2").Value = "mermiden@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:mermiden@gmail.com")

# Missing "@"
$ws.Range("A3").Value = "mermidengmail.com"

# Missing "." before the domain suffix -> still gets auto-linked
$ws.Range("A4").Value = "mermiden@gmailcom"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mermiden@gmailcom")

# Blank / whitespace-only input
$ws.Range("A5").Value = "                   "

# Purely numeric input
$ws.Range("A6").Value = 123456789

# Plain alphabetic text
$ws.Range("A7").Value = "asddfgghl"

# Special characters -> still gets auto-linked
$ws.Range("A8").Value = "!@#`$%^&*()"
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:!@#`$%^&*()")

# "@" replaced with "#" and the "." dropped. The hyperlink was originally
# created against the valid-looking address (so its cached display text
# still reads "mermiden@gmail.com"), then the cell text itself was edited
# to the invalid value afterwards, leaving the link's display cache stale.
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:mermiden@gmail.com", "", "", "mermiden@gmail.com")
$ws.Range("A9").Value = "mermiden#gmail com"

# Widen the column so the email text fits.
$ws.Columns.Item(1).ColumnWidth = 23.6

# Leave the cursor on the last populated cell, and make this new sheet the
# selected / active tab, matching the source workbook.
[void]$ws.Range("A9").Select()
